$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: duplicate of row 5 (Abhi0 / Abhi1 / Abhi2)
$ws.Range("A6").Value = "Abhi0"
$ws.Range("B6").Value = "Abhi1"
$ws.Range("C6").Value = "Abhi2"

# Rows 7-9: new strings Abhi_0 / Abhi_1 / Abhi_2
$ws.Range("A7").Value = "Abhi_0"
$ws.Range("B7").Value = "Abhi_1"
$ws.Range("C7").Value = "Abhi_2"

$ws.Range("A8").Value = "Abhi_0"
$ws.Range("B8").Value = "Abhi_1"
$ws.Range("C8").Value = "Abhi_2"

$ws.Range("A9").Value = "Abhi_0"
$ws.Range("B9").Value = "Abhi_1"
$ws.Range("C9").Value = "Abhi_2"
